# Updates for phs000465, 468, 469 (this workbook: phs000468 Sex-Male / Race-Black or African American)
# Refresh the embedded StatQuery/TabQuery SQL text for the Participants, Files and Samples tabs
# (consent-group-aware joins, extra columns, NULL/"-999" handling, etc.) and keep the Stats
# summary query in sync with the new join shape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2 - StatQuery (counts) for the ParticipantsTab row
$ws.Range("C2").Value = @'
SELECT
  s.Studies,
  p.Participants,
  sa.Samples,
  f.Files
FROM
  ( SELECT COUNT(DISTINCT std.dbgap_accession) AS Studies
    FROM df_study std
    WHERE std.dbgap_accession = 'phs000468'
  ) s
CROSS JOIN
  ( SELECT COUNT(DISTINCT prt.participant_id) AS Participants
    FROM df_study std
    JOIN df_consent_group cg ON std.id = cg."study.id"
    JOIN df_participant prt  ON cg.id  = prt."consent_group.id"
    WHERE std.dbgap_accession = 'phs000468'
      AND prt.sex_at_birth = 'Male'
      AND prt.race = 'Black or African American'
  ) p
CROSS JOIN
  ( SELECT COUNT(DISTINCT smp.sample_id) AS Samples
    FROM df_study std
    JOIN df_consent_group cg ON std.id = cg."study.id"
    JOIN df_participant prt  ON cg.id  = prt."consent_group.id"
    JOIN df_sample smp       ON prt.id = smp."participant.id"
    WHERE std.dbgap_accession = 'phs000468'
      AND prt.sex_at_birth = 'Male'
      AND prt.race = 'Black or African American'
  ) sa
CROSS JOIN
  ( SELECT COUNT(DISTINCT cmf.id) AS Files
    FROM df_clinical_measure_file cmf
    JOIN df_study std2 ON cmf."study.id" = std2.id
    WHERE std2.dbgap_accession = 'phs000468'
  ) f;
'@

# B5 - TabQuery (Files) for the FilesTab row
$ws.Range("B5").Value = @'
WITH file_data AS (
  SELECT
    cmf.file_name,
    REPLACE(COALESCE(cmf.data_category, ''), ';', ', ') AS data_category,
    COALESCE(cmf.file_description, '')                  AS file_description,
    cmf.file_type,
    cmf.file_access,
    COALESCE(cmf.file_size, 0)                          AS file_size,
    cmf."study.id"                                      AS study_row_id,
    cmf."participant.id"                                AS participant_row_id
  FROM df_clinical_measure_file cmf
)
SELECT DISTINCT
  fd.file_name                                 AS "File Name",
  fd.data_category                             AS "Data Category",
  fd.file_description                          AS "File Description",
  fd.file_type                                 AS "File Type",
  CASE
    WHEN fd.file_size >= 1024 * 1024 * 1024 THEN
      CASE
        WHEN ROUND(fd.file_size / (1024.0 * 1024.0 * 1024.0), 2)
             = CAST(ROUND(fd.file_size / (1024.0 * 1024.0 * 1024.0), 0) AS INT)
        THEN CAST(CAST(ROUND(fd.file_size / (1024.0 * 1024.0 * 1024.0), 0) AS INT) AS TEXT) || ' GB'
        ELSE ROUND(fd.file_size / (1024.0 * 1024.0 * 1024.0), 2) || ' GB'
      END
    WHEN fd.file_size >= 1024 * 1024 THEN
      CASE
        WHEN ROUND(fd.file_size / (1024.0 * 1024.0), 2)
             = CAST(ROUND(fd.file_size / (1024.0 * 1024.0), 0) AS INT)
        THEN CAST(CAST(ROUND(fd.file_size / (1024.0 * 1024.0), 0) AS INT) AS TEXT) || ' MB'
        ELSE ROUND(fd.file_size / (1024.0 * 1024.0), 2) || ' MB'
      END
    WHEN fd.file_size >= 1024 THEN
      CASE
        WHEN ROUND(fd.file_size / 1024.0, 2)
             = CAST(ROUND(fd.file_size / 1024.0, 0) AS INT)
        THEN CAST(CAST(ROUND(fd.file_size / 1024.0, 0) AS INT) AS TEXT) || ' KB'
        ELSE ROUND(fd.file_size / 1024.0, 2) || ' KB'
      END
    ELSE
      CASE
        WHEN ROUND(fd.file_size, 2) = CAST(ROUND(fd.file_size, 0) AS INT)
        THEN CAST(CAST(ROUND(fd.file_size, 0) AS INT) AS TEXT) || ' B'
        ELSE ROUND(fd.file_size, 2) || ' B'
      END
  END                                           AS "File Size",
  COALESCE(fd.file_access, '')                  AS "File Access",
  std.dbgap_accession                           AS "Study ID",
  COALESCE(prt.participant_id, '')              AS "Participant ID",
  ''                                            AS "Sample ID"
FROM file_data fd
JOIN df_study std
  ON std.id = fd.study_row_id
LEFT JOIN df_participant prt
  ON prt.id = fd.participant_row_id
WHERE std.dbgap_accession = 'phs000468'
  AND EXISTS (
    SELECT 1
    FROM df_consent_group cg
    JOIN df_participant prt2 ON cg.id = prt2."consent_group.id"
    WHERE cg."study.id" = std.id
      AND prt2.sex_at_birth = 'Male'
      AND prt2.race = 'Black or African American'
  )
ORDER BY fd.file_name ASC;
'@

# B2 - TabQuery (Participants) for the ParticipantsTab row
$ws.Range("B2").Value = @'
WITH diagnosis1 AS (
  SELECT 
    d."participant.id",
    GROUP_CONCAT(
      CASE 
        WHEN COALESCE(CAST(d.age_at_diagnosis AS TEXT), '') IN ('', '-999') THEN NULL
        ELSE CAST(d.age_at_diagnosis AS TEXT)
      END,
      ';'
    ) AS age_days
  FROM df_diagnosis d
  WHERE d."participant.id" IS NOT NULL
  GROUP BY d."participant.id"
),
diagnosis2 AS (
  SELECT "participant.id",
         GROUP_CONCAT(diagnosis, ';') AS diag
  FROM (
    SELECT DISTINCT "participant.id", COALESCE(diagnosis, '') AS diagnosis
    FROM df_diagnosis
    WHERE "participant.id" IS NOT NULL
  )
  GROUP BY "participant.id"
),
diagnosis3 AS (
  SELECT "participant.id",
         GROUP_CONCAT(anatomic_site, ';') AS ant_site
  FROM (
    SELECT DISTINCT "participant.id", COALESCE(anatomic_site, '') AS anatomic_site
    FROM df_diagnosis
    WHERE "participant.id" IS NOT NULL
  )
  GROUP BY "participant.id"
),
diagnosis4 AS (
  SELECT "participant.id",
         GROUP_CONCAT(diagnosis_category, ';') AS diagnosis_category
  FROM (
    SELECT DISTINCT "participant.id", COALESCE(diagnosis_category, '') AS diagnosis_category
    FROM df_diagnosis
    WHERE "participant.id" IS NOT NULL
  )
  GROUP BY "participant.id"
),
survival_agg AS (
  SELECT s."participant.id",
         MAX(CASE WHEN COALESCE(s.last_known_survival_status, '') <> '' 
                  THEN s.last_known_survival_status ELSE '' END) AS last_known_survival_status
  FROM df_survival s
  GROUP BY s."participant.id"
)
SELECT DISTINCT
  prt.participant_id                          AS "Participant ID",
  std.dbgap_accession                         AS "Study ID",
  COALESCE(prt.sex_at_birth, '')              AS "Sex",
  COALESCE(prt.race, '')                      AS "Race",
  dgn2.diag                                   AS "Diagnosis",
  dgn3.ant_site                               AS "Diagnosis Anatomic Site",
  dgn4.diagnosis_category                     AS "Diagnosis Category",
  COALESCE(
    CASE 
      WHEN COALESCE(dgn1.age_days, '') = '' THEN 'Not Reported'
      ELSE dgn1.age_days
    END,
    'Not Reported'
  )                                            AS "Age at Diagnosis (days)",
  NULL                                         AS "Treatment Type",
  COALESCE(sva.last_known_survival_status, '') AS "Last Known Survival Status"
FROM df_study std
JOIN df_consent_group cg
  ON std.id = cg."study.id"
JOIN df_participant prt
  ON cg.id = prt."consent_group.id"
LEFT JOIN diagnosis1     dgn1 ON prt.id = dgn1."participant.id"
LEFT JOIN diagnosis2     dgn2 ON prt.id = dgn2."participant.id"
LEFT JOIN diagnosis3     dgn3 ON prt.id = dgn3."participant.id"
LEFT JOIN diagnosis4     dgn4 ON prt.id = dgn4."participant.id"
LEFT JOIN survival_agg   sva  ON prt.id = sva."participant.id"
WHERE std.dbgap_accession = 'phs000468'
  AND prt.race = 'Black or African American'
  AND prt.sex_at_birth = 'Male'
ORDER BY prt.participant_id ASC;
'@

# B4 - TabQuery (Samples) for the SamplesTab row
$ws.Range("B4").Value = @'
WITH diag_by_sample AS (
  SELECT 
    d.sample_row_id,
    GROUP_CONCAT(d.diagnosis, '; ') AS sample_diagnosis
  FROM (
    SELECT DISTINCT 
      dd."sample.id" AS sample_row_id,
      COALESCE(dd.diagnosis, '') AS diagnosis
    FROM df_diagnosis dd
    WHERE dd."sample.id" IS NOT NULL
  ) d
  GROUP BY d.sample_row_id
),
diagcat_by_sample AS (
  SELECT
    c.sample_row_id,
    GROUP_CONCAT(c.diagnosis_category, '; ') AS sample_diag_category
  FROM (
    SELECT DISTINCT
      dd."sample.id" AS sample_row_id,
      COALESCE(dd.diagnosis_category, '') AS diagnosis_category
    FROM df_diagnosis dd
    WHERE dd."sample.id" IS NOT NULL
  ) c
  GROUP BY c.sample_row_id
)
SELECT DISTINCT
  smp.sample_id                               AS "Sample ID",
  prt.participant_id                           AS "Participant ID",
  std.dbgap_accession                          AS "Study ID",
  COALESCE(smp.anatomic_site, '')              AS "Sample Anatomic Site",
  COALESCE(
    CASE 
      WHEN smp.participant_age_at_collection = -999 OR smp.participant_age_at_collection IS NULL 
        THEN 'Not Reported'
      ELSE CAST(smp.participant_age_at_collection AS TEXT)
    END, 
    'Not Reported'
  )                                            AS "Age at Sample Collection (days)",
  COALESCE(smp.sample_tumor_status, '')        AS "Sample Tumor Status",
  COALESCE(smp.tumor_classification, '')       AS "Sample Tumor Classification",
  COALESCE(dbs.sample_diagnosis, '')           AS "Sample Diagnosis",
  COALESCE(dcs.sample_diag_category, '')       AS "Diagnosis Category"
FROM df_study std
JOIN df_consent_group cg
  ON std.id = cg."study.id"
JOIN df_participant prt
  ON cg.id = prt."consent_group.id"
LEFT JOIN df_sample smp
  ON prt.id = smp."participant.id"
LEFT JOIN diag_by_sample dbs
  ON smp.id = dbs.sample_row_id
LEFT JOIN diagcat_by_sample dcs
  ON smp.id = dcs.sample_row_id
WHERE std.dbgap_accession = 'phs000468'
  AND prt.race = 'Black or African American'
  AND prt.sex_at_birth = 'Male'
  AND smp.sample_id IS NOT NULL
ORDER BY smp.sample_id ASC;
'@

# Row heights settle back to Excel's auto-fit result (these long, multi-line
# SQL blocks push rows 3-5 up against Excel's ~409pt row-height ceiling, and
# row 2 - shared with the equally long StatQuery text in C2 - sits just under it)
$ws.Rows.Item(2).RowHeight = 355
$ws.Rows.Item(3).RowHeight = 409.6
$ws.Rows.Item(4).RowHeight = 409.6
$ws.Rows.Item(5).RowHeight = 409.6

# Update the selection/active view to match the saved workbook state (row 5, B5)
$ws.Range("B5").Select()
